$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H5").Value = 69.416664
$ws1.Range("I5").Value = 66.63636
$ws1.Range("K5").Value = 66.63636
$ws1.Range("M5").Value = 48.36364
$ws1.Range("H18").Value = 1056.6666
$ws1.Range("I18").Value = 1056.6666
$ws1.Range("K18").Value = 1056.6666
$ws1.Range("M18").Value = -772.6666
$ws1.Range("H26").Value = 4000
$ws1.Range("J26").Value = 0
$ws1.Range("L26").Value = 0
$ws1.Range("H92").Value = 1783.3
$ws1.Range("I92").Value = 1783.3
$ws1.Range("J92").Value = 0
$ws1.Range("K92").Value = 1783.3
$ws1.Range("L92").Value = 0
$ws1.Range("N92").Value = -535.3
$ws1.Range("H98").Value = 1173
$ws1.Range("I98").Value = 1173
$ws1.Range("K98").Value = 1173
$ws1.Range("M98").Value = 325
$ws1.Range("H122").Value = 1173
$ws1.Range("I122").Value = 1173
$ws1.Range("K122").Value = 3519
$ws1.Range("M122").Value = -1069
$ws1.Range("H135").Value = 1660
$ws1.Range("I135").Value = 1660
$ws1.Range("K135").Value = 14940
$ws1.Range("M135").Value = -12405
$ws1.Range("H137").Value = 3147.625
$ws1.Range("I137").Value = 1791.6666
$ws1.Range("J137").Value = 3460.5386
$ws1.Range("K137").Value = 5374.9998
$ws1.Range("L137").Value = 10381.6158
$ws1.Range("M137").Value = -2824.9998
$ws1.Range("N137").Value = -15481.6158
$ws1.Range("N26").ClearContents()
$ws1.Range("M92").ClearContents()

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H61").Value = 6227.273
$ws2.Range("I61").Value = 5000
$ws2.Range("K61").Value = 5000
$ws2.Range("M61").Value = -4788
$ws2.Range("H68").Value = 70000
$ws2.Range("J68").Value = 70000
$ws2.Range("L68").Value = 70000
$ws2.Range("N68").Value = -71622
$ws2.Range("H71").Value = 70000
$ws2.Range("J71").Value = 70000
$ws2.Range("L71").Value = 210000
$ws2.Range("N71").Value = -218112
$ws2.Range("H80").Value = 80000
$ws2.Range("I80").Value = 40000
$ws2.Range("J80").Value = 100000
$ws2.Range("K80").Value = 40000
$ws2.Range("L80").Value = 100000
$ws2.Range("M80").Value = -39002
$ws2.Range("N80").Value = -101996
$ws2.Range("H83").Value = 80000
$ws2.Range("I83").Value = 40000
$ws2.Range("J83").Value = 100000
$ws2.Range("K83").Value = 120000
$ws2.Range("L83").Value = 300000
$ws2.Range("M83").Value = -115008
$ws2.Range("N83").Value = -309984
$ws2.Range("H97").Value = 1528.5555
$ws2.Range("I97").Value = 1597.875
$ws2.Range("K97").Value = 1597.875
$ws2.Range("M97").Value = -1101.875
$ws2.Range("H136").Value = 6227.273
$ws2.Range("I136").Value = 5000
$ws2.Range("K136").Value = 15000
$ws2.Range("M136").Value = -12450

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H86").Value = 629.375
$ws3.Range("I86").Value = 597.8570999999999
$ws3.Range("K86").Value = 597.8570999999999
$ws3.Range("M86").Value = 525.1429000000001
$ws3.Range("H89").Value = 629.375
$ws3.Range("I89").Value = 597.8570999999999
$ws3.Range("K89").Value = 2989.2855
$ws3.Range("M89").Value = 2626.7145
$ws3.Range("H105").Value = 1758.3334
$ws3.Range("I105").Value = 1712.5
$ws3.Range("J105").Value = 1850
$ws3.Range("K105").Value = 1712.5
$ws3.Range("L105").Value = 1850
$ws3.Range("M105").Value = 34.5
$ws3.Range("N105").Value = -5344

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H31").Value = 4297.1035
$ws4.Range("I31").Value = 2054.5334
$ws4.Range("J31").Value = 5079.3955
$ws4.Range("K31").Value = 2054.5334
$ws4.Range("L31").Value = 5079.3955
$ws4.Range("M31").Value = -1759.5334
$ws4.Range("N31").Value = -5669.3955
$ws4.Range("H34").Value = 4297.1035
$ws4.Range("I34").Value = 2054.5334
$ws4.Range("J34").Value = 5079.3955
$ws4.Range("K34").Value = 2054.5334
$ws4.Range("L34").Value = 5079.3955
$ws4.Range("M34").Value = -1852.5334
$ws4.Range("N34").Value = -5483.3955
$ws4.Range("H107").Value = 450.16666
$ws4.Range("I107").Value = 375.55554
$ws4.Range("K107").Value = 375.55554
$ws4.Range("M107").Value = 1544.44446
$ws4.Range("H132").Value = 1312.5
$ws4.Range("I132").Value = 1312.5
$ws4.Range("K132").Value = 3937.5
$ws4.Range("M132").Value = -1407.5

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H68").Value = 3581.3333
$ws5.Range("I68").Value = 2247.5
$ws5.Range("K68").Value = 6742.5
$ws5.Range("M68").Value = -5931.5
$ws5.Range("H71").Value = 3581.3333
$ws5.Range("I71").Value = 2247.5
$ws5.Range("K71").Value = 20227.5
$ws5.Range("M71").Value = -16171.5

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H132").Value = 6285.143
$ws6.Range("J132").Value = 2997
$ws6.Range("L132").Value = 8991
$ws6.Range("N132").Value = -14051

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H7").Value = 6964.7144
$ws7.Range("I7").Value = 6499.6
$ws7.Range("J7").Value = 8127.5
$ws7.Range("K7").Value = 6499.6
$ws7.Range("L7").Value = 8127.5
$ws7.Range("M7").Value = -6387.6
$ws7.Range("N7").Value = -8351.5
$ws7.Range("H22").Value = 3605.111
$ws7.Range("I22").Value = 600
$ws7.Range("J22").Value = 4463.7144
$ws7.Range("K22").Value = 600
$ws7.Range("L22").Value = 4463.7144
$ws7.Range("M22").Value = -305
$ws7.Range("N22").Value = -5053.7144
$ws7.Range("H27").Value = 3605.111
$ws7.Range("I27").Value = 600
$ws7.Range("J27").Value = 4463.7144
$ws7.Range("K27").Value = 600
$ws7.Range("L27").Value = 4463.7144
$ws7.Range("M27").Value = -493
$ws7.Range("N27").Value = -4677.7144
$ws7.Range("H40").Value = 1999
$ws7.Range("I40").Value = 1999
$ws7.Range("K40").Value = 1999
$ws7.Range("M40").Value = -1863
$ws7.Range("H46").Value = 5369.4287
$ws7.Range("I46").Value = 4632
$ws7.Range("J46").Value = 6106.857
$ws7.Range("K46").Value = 4632
$ws7.Range("L46").Value = 6106.857
$ws7.Range("M46").Value = -4444
$ws7.Range("N46").Value = -6482.857
$ws7.Range("H68").Value = 8220.777
$ws7.Range("I68").Value = 7332.3335
$ws7.Range("J68").Value = 8665
$ws7.Range("K68").Value = 7332.3335
$ws7.Range("L68").Value = 8665
$ws7.Range("M68").Value = -6583.3335
$ws7.Range("N68").Value = -10163
$ws7.Range("H71").Value = 8220.777
$ws7.Range("I71").Value = 7332.3335
$ws7.Range("J71").Value = 8665
$ws7.Range("K71").Value = 36661.6675
$ws7.Range("L71").Value = 43325
$ws7.Range("M71").Value = -32917.6675
$ws7.Range("N71").Value = -50813
$ws7.Range("H82").Value = 3500
$ws7.Range("I82").Value = 3500
$ws7.Range("J82").Value = 0
$ws7.Range("K82").Value = 3500
$ws7.Range("L82").Value = 0
$ws7.Range("N82").Value = -3139
$ws7.Range("H85").Value = 3500
$ws7.Range("I85").Value = 3500
$ws7.Range("J85").Value = 0
$ws7.Range("K85").Value = 3500
$ws7.Range("L85").Value = 0
$ws7.Range("N85").Value = -2252
$ws7.Range("H126").Value = 6964.7144
$ws7.Range("I126").Value = 6499.6
$ws7.Range("J126").Value = 8127.5
$ws7.Range("K126").Value = 19498.8
$ws7.Range("L126").Value = 24382.5
$ws7.Range("M126").Value = -17028.8
$ws7.Range("N126").Value = -29322.5
$ws7.Range("H136").Value = 1947.625
$ws7.Range("I136").Value = 1863.5555
$ws7.Range("J136").Value = 2199.8333
$ws7.Range("K136").Value = 5590.666499999999
$ws7.Range("L136").Value = 6599.499899999999
$ws7.Range("M136").Value = -3040.666499999999
$ws7.Range("N136").Value = -11699.4999
$ws7.Range("M82").ClearContents()
$ws7.Range("M85").ClearContents()

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H40").Value = 60000
$ws8.Range("I40").Value = 60000
$ws8.Range("K40").Value = 60000
$ws8.Range("M40").Value = -59851
$ws8.Range("H82").Value = 36250
$ws8.Range("J82").Value = 36250
$ws8.Range("L82").Value = 36250
$ws8.Range("N82").Value = -37016
$ws8.Range("H85").Value = 36250
$ws8.Range("J85").Value = 36250
$ws8.Range("L85").Value = 36250
$ws8.Range("N85").Value = -38902
$ws8.Range("H122").Value = 2245.6667
$ws8.Range("I122").Value = 2166.2632
$ws8.Range("K122").Value = 6498.7896
$ws8.Range("M122").Value = -4048.7896
$ws8.Range("H126").Value = 4346.1577
$ws8.Range("I126").Value = 2752.4546
$ws8.Range("J126").Value = 6537.5
$ws8.Range("K126").Value = 8257.363799999999
$ws8.Range("L126").Value = 19612.5
$ws8.Range("M126").Value = -5787.363799999999
$ws8.Range("N126").Value = -24552.5
$ws8.Range("H132").Value = 1660.75
$ws8.Range("I132").Value = 1660.75
$ws8.Range("K132").Value = 4982.25
$ws8.Range("M132").Value = -2452.25
$ws8.Range("H136").Value = 3801.3333
$ws8.Range("I136").Value = 2093.2632
$ws8.Range("J136").Value = 6119.4287
$ws8.Range("K136").Value = 6279.7896
$ws8.Range("L136").Value = 18358.2861
$ws8.Range("M136").Value = -3729.7896
$ws8.Range("N136").Value = -23458.2861
$ws8.Range("H139").Value = 81399.5
$ws8.Range("J139").Value = 81399.5
$ws8.Range("L139").Value = 81399.5
$ws8.Range("N139").Value = -91679.5

Write-Host "Applied all Halicarnassus_Profits updates."
